# Scheduled Typhon_Profits market-price refresh: update cached leve profit figures
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6259398
$ws.Range("J17").Value = 6676564.5
$ws.Range("L17").Value = 20029693.5
$ws.Range("N17").Value = -20030029.5
$ws.Range("H40").Value = 1387.0588
$ws.Range("I40").Value = 776.2
$ws.Range("K40").Value = 776.2
$ws.Range("M40").Value = -601.2
$ws.Range("H62").Value = 4266.4165
$ws.Range("I62").Value = 3163.125
$ws.Range("J62").Value = 6473
$ws.Range("K62").Value = 3163.125
$ws.Range("L62").Value = 6473
$ws.Range("M62").Value = -2539.125
$ws.Range("N62").Value = -7721
$ws.Range("H65").Value = 4266.4165
$ws.Range("I65").Value = 3163.125
$ws.Range("J65").Value = 6473
$ws.Range("K65").Value = 15815.625
$ws.Range("L65").Value = 32365
$ws.Range("M65").Value = -12695.625
$ws.Range("N65").Value = -38605
$ws.Range("H69").Value = 1515.5
$ws.Range("I69").Value = 1200
$ws.Range("J69").Value = 1526.3793
$ws.Range("K69").Value = 3600
$ws.Range("L69").Value = 4579.1379
$ws.Range("M69").Value = -2726
$ws.Range("N69").Value = -6327.1379
$ws.Range("H72").Value = 1515.5
$ws.Range("I72").Value = 1200
$ws.Range("J72").Value = 1526.3793
$ws.Range("K72").Value = 10800
$ws.Range("L72").Value = 13737.4137
$ws.Range("M72").Value = -6432
$ws.Range("N72").Value = -22473.4137
$ws.Range("H98").Value = 1257
$ws.Range("I98").Value = 1071.25
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1071.25
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 426.75
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 1257
$ws.Range("I122").Value = 1071.25
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3213.75
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -763.75
$ws.Range("N122").Value = -10900
$ws.Range("H129").Value = 137904
$ws.Range("J129").Value = 162327.84
$ws.Range("L129").Value = 486983.52
$ws.Range("N129").Value = -496983.52
$ws.Range("H137").Value = 1867.0526
$ws.Range("I137").Value = 1594.1
$ws.Range("J137").Value = 2170.3333
$ws.Range("K137").Value = 4782.299999999999
$ws.Range("L137").Value = 6510.999899999999
$ws.Range("M137").Value = -2232.299999999999
$ws.Range("N137").Value = -11610.9999
$ws.Range("H139").Value = 50513.332
$ws.Range("J139").Value = 50513.332
$ws.Range("L139").Value = 50513.332
$ws.Range("N139").Value = -60793.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1894.6389
$ws.Range("I86").Value = 1677.5
$ws.Range("J86").Value = 2235.8572
$ws.Range("K86").Value = 1677.5
$ws.Range("L86").Value = 2235.8572
$ws.Range("M86").Value = -554.5
$ws.Range("N86").Value = -4481.8572
$ws.Range("H89").Value = 1894.6389
$ws.Range("I89").Value = 1677.5
$ws.Range("J89").Value = 2235.8572
$ws.Range("K89").Value = 8387.5
$ws.Range("L89").Value = 11179.286
$ws.Range("M89").Value = -2771.5
$ws.Range("N89").Value = -22411.286
$ws.Range("H105").Value = 2274565.8
$ws.Range("I105").Value = 1708.1666
$ws.Range("K105").Value = 1708.1666
$ws.Range("M105").Value = 38.83339999999998
$ws.Range("H134").Value = 4650.846
$ws.Range("I134").Value = 4963.4165
$ws.Range("K134").Value = 14890.2495
$ws.Range("M134").Value = -12355.2495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3045.1875
$ws.Range("I31").Value = 1532.5
$ws.Range("K31").Value = 1532.5
$ws.Range("M31").Value = -1237.5
$ws.Range("H34").Value = 3045.1875
$ws.Range("I34").Value = 1532.5
$ws.Range("K34").Value = 1532.5
$ws.Range("M34").Value = -1330.5
$ws.Range("H99").Value = 3096.889
$ws.Range("J99").Value = 4628.5713
$ws.Range("L99").Value = 4628.5713
$ws.Range("N99").Value = -7624.5713
$ws.Range("H122").Value = 1590.2858
$ws.Range("I122").Value = 1538
$ws.Range("K122").Value = 4614
$ws.Range("M122").Value = -2164
$ws.Range("H126").Value = 3096.889
$ws.Range("J126").Value = 4628.5713
$ws.Range("L126").Value = 13885.7139
$ws.Range("N126").Value = -18825.7139
$ws.Range("H132").Value = 2480.2258
$ws.Range("I132").Value = 1867.36
$ws.Range("J132").Value = 5033.8335
$ws.Range("K132").Value = 5602.08
$ws.Range("L132").Value = 15101.5005
$ws.Range("M132").Value = -3072.08
$ws.Range("N132").Value = -20161.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12521.375
$ws.Range("I2").Value = 12521.375
$ws.Range("K2").Value = 75128.25
$ws.Range("M2").Value = -75015.25
$ws.Range("H5").Value = 1320.6586
$ws.Range("I5").Value = 949.9259
$ws.Range("J5").Value = 2035.6428
$ws.Range("K5").Value = 2849.7777
$ws.Range("L5").Value = 6106.928400000001
$ws.Range("M5").Value = -2737.7777
$ws.Range("N5").Value = -6330.928400000001
$ws.Range("H68").Value = 744.5
$ws.Range("J68").Value = 490
$ws.Range("L68").Value = 1470
$ws.Range("N68").Value = -3092
$ws.Range("H71").Value = 744.5
$ws.Range("J71").Value = 490
$ws.Range("L71").Value = 4410
$ws.Range("N71").Value = -12522
$ws.Range("H131").Value = 701.61
$ws.Range("J131").Value = 724.5761
$ws.Range("L131").Value = 2173.7283
$ws.Range("N131").Value = -12253.7283
$ws.Range("H134").Value = 3039.3076
$ws.Range("I134").Value = 1887.9375
$ws.Range("J134").Value = 4881.5
$ws.Range("K134").Value = 5663.8125
$ws.Range("L134").Value = 14644.5
$ws.Range("M134").Value = -593.8125
$ws.Range("N134").Value = -24784.5
$ws.Range("H135").Value = 1320.6586
$ws.Range("I135").Value = 949.9259
$ws.Range("J135").Value = 2035.6428
$ws.Range("K135").Value = 8549.3331
$ws.Range("L135").Value = 18320.7852
$ws.Range("M135").Value = -6014.3331
$ws.Range("N135").Value = -23390.7852
$ws.Range("H136").Value = 3642.25
$ws.Range("I136").Value = 932.5
$ws.Range("J136").Value = 4997.125
$ws.Range("K136").Value = 2797.5
$ws.Range("L136").Value = 14991.375
$ws.Range("M136").Value = 2302.5
$ws.Range("N136").Value = -25191.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 570.5714
$ws.Range("I107").Value = 570.5714
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 570.5714
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1349.4286
$ws.Range("H109").Value = 28935
$ws.Range("J109").Value = 28935
$ws.Range("L109").Value = 28935
$ws.Range("N107").ClearContents()
$ws.Range("N109").Value = -31015

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2555.7437
$ws.Range("I40").Value = 2341.4194
$ws.Range("J40").Value = 3386.25
$ws.Range("K40").Value = 2341.4194
$ws.Range("L40").Value = 3386.25
$ws.Range("M40").Value = -2205.4194
$ws.Range("N40").Value = -3658.25
$ws.Range("H93").Value = 1535.6364
$ws.Range("I93").Value = 1715.3334
$ws.Range("K93").Value = 1715.3334
$ws.Range("M93").Value = -467.3334
